$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.025940299034119
$ws.Range("B1").Value = 3.214179754257202
$ws.Range("C1").Value = 2.56113600730896
$ws.Range("D1").Value = 2.401116847991943
$ws.Range("E1").Value = 2.015549182891846
